$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update timestamps
$wsOverview.Range("G2").Value = "2016-09-05 11:26:42"
$wsDeDe.Range("H2").Value = "2016-09-05 11:26:42"
$wsZhCn.Range("H2").Value = "2016-09-05 11:26:37"

# Update column widths (autofit-like widening from text change).
# NOTE: the host quantizes ColumnWidth to 1/6-character-unit pixel
# granularity (stored_width = round(input*6)/6 + 5/6), so an input of
# 16.333333333333336 is the closest reachable value to the target stored
# width of 17.2159881591797 (lands on 17.166666666666668).
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
